$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-order the three matches currently sitting in rows 97-99 -------
# Old layout:
#   row 97 -> Ankaragucu  vs Samsunspor
#   row 98 -> Antalyaspor vs Basaksehir
#   row 99 -> Pendikspor  vs Fenerbahce
# New layout:
#   row 97 -> Antalyaspor vs Basaksehir   (was row 98)
#   row 98 -> Pendikspor  vs Fenerbahce   (was row 99)
#   row 99 -> Ankaragucu  vs Samsunspor   (was row 97)
# Columns A:E (index/country/tournament/season/kickoff date) stay put for
# each row position - only the match data in F:V rotates between rows.

$row97 = $ws.Range("F97:V97").Value()
$row98 = $ws.Range("F98:V98").Value()
$row99 = $ws.Range("F99:V99").Value()

$ws.Range("F99:V99").Value = $row97
$ws.Range("F97:V97").Value = $row98
$ws.Range("F98:V98").Value = $row99

# --- 2. Append the new match as row 101 -----------------------------------
# Copy row 100 first so the new row inherits the correct cell styles
# (bold/bordered index in column A, date number-format in column E).
$ws.Range("A100:V100").Copy($ws.Range("A101:V101"))

$ws.Cells.Item(101, 1).Value = 100
$ws.Cells.Item(101, 2).Value = "turkey"
$ws.Cells.Item(101, 3).Value = "super-lig"
$ws.Cells.Item(101, 4).Value = "2023-2024"
$ws.Cells.Item(101, 5).Value = 45233.75
$ws.Cells.Item(101, 6).Value = "Galatasaray"
$ws.Cells.Item(101, 7).Value = 2
$ws.Cells.Item(101, 8).Value = "Kasimpasa"
$ws.Cells.Item(101, 9).Value = 1
$ws.Cells.Item(101, 10).Value = 1.24
$ws.Cells.Item(101, 11).Value = "28/10/2023 18:13"
$ws.Cells.Item(101, 12).Value = 1.21
$ws.Cells.Item(101, 13).Value = "03/11/2023 17:58"
$ws.Cells.Item(101, 14).Value = 7.31
$ws.Cells.Item(101, 15).Value = "28/10/2023 18:13"
$ws.Cells.Item(101, 16).Value = 7.69
$ws.Cells.Item(101, 17).Value = "03/11/2023 17:58"
$ws.Cells.Item(101, 18).Value = 11.06
$ws.Cells.Item(101, 19).Value = "28/10/2023 18:13"
$ws.Cells.Item(101, 20).Value = 12.12
$ws.Cells.Item(101, 21).Value = "03/11/2023 17:58"
$ws.Cells.Item(101, 22).Value = "https://www.betexplorer.com/football/turkey/super-lig/galatasaray-kasimpasa/hz3dxW51/"
